$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new NARA columns L..P
$ws.Range("L1").Value = "NARA_Format Name"
$ws.Range("M1").Value = "NARA_PRONOM URL"
$ws.Range("N1").Value = "NARA_Risk Level"
$ws.Range("O1").Value = "NARA_Proposed Preservation Plan"
$ws.Range("P1").Value = "NARA_Match_Type"

# Row 2
$ws.Range("L2").Value = "Exchangeable Image File Format Compressed 2.1"
$ws.Range("M2").Value = "https://www.nationalarchives.gov.uk/pronom/x-fmt/390"
$ws.Range("N2").Value = "Low Risk"
$ws.Range("O2").Value = "Retain"
$ws.Range("P2").Value = "PRONOM"

# Row 3
$ws.Range("L3").Value = "Exchangeable Image File Format Compressed 2.1"
$ws.Range("M3").Value = "https://www.nationalarchives.gov.uk/pronom/x-fmt/390"
$ws.Range("N3").Value = "Low Risk"
$ws.Range("O3").Value = "Retain"
$ws.Range("P3").Value = "PRONOM"

# Row 4
$ws.Range("L4").Value = "JPEG File Interchange Format 1.01"
$ws.Range("M4").Value = "https://www.nationalarchives.gov.uk/pronom/fmt/43"
$ws.Range("N4").Value = "Low Risk"
$ws.Range("O4").Value = "Retain"
$ws.Range("P4").Value = "PRONOM and Version"

# Selection update to match target (active cell P1 selected)
$ws.Range("P1").Select()
